$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Compass Error data (was row 5)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Compass Error Compass Error Compass data error Please contact DJI Support ."
$ws.Range("C2").Value = "Please contact DJI Support"
$ws.Range("D2").Value = "7-10"
$ws.Range("E2").Value = "NonEvent"
$ws.Range("F2").Value = "NonEvent"

# Row 3: High wind velocity first sentence (was row 2)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP ."
$ws.Range("C3").Value = "High wind velocity"
$ws.Range("D3").Value = "0-2"
$ws.Range("E3").Value = "Event"
$ws.Range("F3").Value = "Event"

# Row 4: High wind velocity second sentence (was row 3)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP ."
$ws.Range("C4").Value = "Ensure the aircraft remains within your line of sight and fly with caution"
$ws.Range("D4").Value = "3-15"
$ws.Range("E4").Value = "NonEvent"
$ws.Range("F4").Value = "NonEvent"

# Row 5: High wind velocity third sentence (was row 4)
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP ."
$ws.Range("C5").Value = "Fly with caution and land in a safe place ASAP"
$ws.Range("D5").Value = "19-28"
$ws.Range("E5").Value = "NonEvent"
$ws.Range("F5").Value = "NonEvent"
